$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''69.660.19'
$ws.Range("E2").Value = '''  +1.90%  '
$ws.Range("D3").Value = '''3.888.31'
$ws.Range("E3").Value = '''  +1.46%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '''  +0.03%  '
$ws.Range("D5").Value = '''604.40'
$ws.Range("E5").Value = '''  +0.89%  '
$ws.Range("D6").Value = '''170.26'
$ws.Range("E6").Value = '''  +5.31%  '
$ws.Range("D7").Value = '''3.887.77'
$ws.Range("E7").Value = '''  +1.44%  '
$ws.Range("E8").Value = '''  +0.13%  '
$ws.Range("E9").Value = '''  +1.34%  '
$ws.Range("E10").Value = '''  +1.69%  '
$ws.Range("E11").Value = '''  +1.37%  '
$ws.Range("D12").Value = '''0.466'
$ws.Range("E12").Value = '''  +2.04%  '
$ws.Range("D13").Value = '''0.0000255'
$ws.Range("E13").Value = '''  +5.23%  '
$ws.Range("D14").Value = '''38.13'
$ws.Range("E14").Value = '''  +4.12%  '
$ws.Range("D15").Value = '''4.548.74'
$ws.Range("E15").Value = '''  +1.63%  '
$ws.Range("D16").Value = '''3.895.67'
$ws.Range("E16").Value = '''  +1.48%  '
$ws.Range("D17").Value = '''69.702.07'
$ws.Range("E17").Value = '''  +1.60%  '
$ws.Range("D18").Value = '''18.68'
$ws.Range("E18").Value = '''  +9.44%  '
$ws.Range("D19").Value = '''7.63'
$ws.Range("E19").Value = '''  +1.71%  '
$ws.Range("E20").Value = '''  -0.71%  '
$ws.Range("D21").Value = '''11.16'
$ws.Range("E21").Value = '''  -0.90%  '
$ws.Range("D22").Value = '''489.37'
$ws.Range("E22").Value = '''  +1.22%  '
$ws.Range("D23").Value = '''0.747'
$ws.Range("E23").Value = '''  +4.75%  '
$ws.Range("D24").Value = '''0.0000165'
$ws.Range("E24").Value = '''  +3.36%  '
$ws.Range("D25").Value = '''85.23'
$ws.Range("E25").Value = '''  +1.75%  '
$ws.Range("D26").Value = '''2.31'
$ws.Range("E26").Value = '''  +4.12%  '
$ws.Range("D27").Value = '''12.35'
$ws.Range("E27").Value = '''  +2.46%  '
$ws.Range("D28").Value = '''10.11'
$ws.Range("E28").Value = '''  +2.28%  '
$ws.Range("E29").Value = '''  +0.28%  '
$ws.Range("D30").Value = '''2.97'
$ws.Range("E30").Value = '''  +1.34%  '
$ws.Range("D31").Value = '''4.043.69'
$ws.Range("E31").Value = '''  +1.34%  '
$ws.Range("D32").Value = '''2.41'
$ws.Range("E32").Value = '''  +2.53%  '
$ws.Range("E33").Value = '''  +0.35%  '
$ws.Range("D34").Value = '''31.83'
$ws.Range("E34").Value = '''  -0.21%  '
$ws.Range("D35").Value = '''3.859.58'
$ws.Range("E35").Value = '''  +2.03%  '
$ws.Range("E36").Value = '''  +0.75%  '
$ws.Range("D37").Value = '''6.11'
$ws.Range("E37").Value = '''  +4.68%  '
$ws.Range("D38").Value = '''1.03'
$ws.Range("E38").Value = '''  +0.73%  '
$ws.Range("E39").Value = '''  +2.29%  '
$ws.Range("D40").Value = '''3.35'
$ws.Range("E40").Value = '''  +14.33%  '
$ws.Range("D41").Value = '''0.999'
$ws.Range("E41").Value = '''  -0.09%  '
$ws.Range("D42").Value = '''0.328'
$ws.Range("E42").Value = '''  +3.82%  '
$ws.Range("D43").Value = '''2.08'
$ws.Range("E43").Value = '''  +6.28%  '
$ws.Range("D44").Value = '''437.36'
$ws.Range("E44").Value = '''  +2.50%  '
$ws.Range("D45").Value = '''48.27'
$ws.Range("E45").Value = '''  -0.41%  '
$ws.Range("D46").Value = '''8.68'
$ws.Range("E46").Value = '''  +4.25%  '
$ws.Range("D48").Value = '''0.000277'
$ws.Range("E48").Value = '''  +23.71%  '
$ws.Range("E49").Value = '''  +3.79%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '''144.15'
$ws.Range("E50").Value = '''  +0.91%  '
$ws.Range("B51").Value = 'Arweave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D51").Value = '''40.39'
$ws.Range("E51").Value = '''  +5.08%  '
